# Deploying to gh-pages from  @ fa727193d3653d1e75f6d95b559963a078b26bde
#
# Adds a new "2023" data column (column T) to the "1.5.1 Number of deaths
# attributed to disasters" worksheet, mirroring the existing column-S
# ("2022") layout/formatting, and tidies up column A:C widths / the saved
# selection that the original workbook had accumulated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Copy the formatting of column S (the most recent existing year)
#    onto the new column T before writing any values, so every T cell
#    inherits the same borders/fonts/number formatting as its S sibling.
# ---------------------------------------------------------------------
$srcFormat = $ws.Range("S3:S34")
$dstFormat = $ws.Range("T3:T34")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Header cell: T4 = 2023 (new year column header)
# ---------------------------------------------------------------------
$ws.Range("T4").Value = 2023

# ---------------------------------------------------------------------
# 3) Data cells T5:T34 — values taken from the 2023 column being added.
#    A "-" string denotes "no data" (matches the other no-data cells in
#    the sheet, e.g. S17, S18, S19, ...).
# ---------------------------------------------------------------------
$values = @{
    5  = 44
    6  = 24
    7  = 20
    8  = "-"
    9  = "-"
    10 = "-"
    11 = 5
    12 = 1
    13 = 4
    14 = 8
    15 = 6
    16 = 2
    17 = 5
    18 = 1
    19 = 4
    20 = 7
    21 = 5
    22 = 2
    23 = "-"
    24 = "-"
    25 = "-"
    26 = 18
    27 = 10
    28 = 8
    29 = "-"
    30 = "-"
    31 = "-"
    32 = 1
    33 = 1
    34 = "-"
}

foreach ($r in $values.Keys | Sort-Object) {
    $ws.Cells.Item($r, 20).Value = $values[$r]
}

# ---------------------------------------------------------------------
# 4) Cosmetic clean-up that came with the same commit: the three label
#    columns (A:C) are unified to one common width, and the leftover
#    "T24" selection that had been saved with the workbook is cleared
#    back to the top-left cell.
# ---------------------------------------------------------------------
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 35

$ws.Range("A1").Select() | Out-Null
